$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dashboard page creation: populate task row 8 (Dashboard task)
$ws.Range("C8").Value = "Jakub Ivan Vanko"
$ws.Range("D8").Value = "UI - Home"
$ws.Range("E8").Value = "Dashboard - overview and charts"

# Move the active selection to E8 (matches author's last-edited cell)
$ws.Range("E8").Select()
